# "se sube nueva data para regresion en preProd"
# Replace the RPM009 regression-test row with a new RGA009 row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# NroCuenta (E2): new account number for this regression run
$ws.Range("E2").Value = 6759658789

# FechaInicio (K2): keep it text (cell is quote-prefixed in the source file)
# so Excel doesn't reinterpret the dd/mm/yyyy string as a date serial.
$ws.Range("K2").Value = "'07/04/2021"

# Patente (Z2): new plate id. The source cell has no explicit cell style,
# so reset formatting to Normal before writing the value.
$ws.Range("Z2").Style = "Normal"
$ws.Range("Z2").Value = "RGA009"

# Motor / Chasis (AA2, AB2): same new identifier, mirrored in both columns
$ws.Range("AA2").Value = "1234567RGA009"
$ws.Range("AB2").Value = "1234567RGA009"

# Move the active selection to E2 (NroCuenta) and scroll the sheet so
# column E is the leftmost visible column, matching the saved view.
$ws.Activate()
$ws.Range("E2").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1

# Restore the (maximized) window geometry recorded with this save.
$excel.Left = -120
$excel.Top = -120
$excel.Width = 20730
$excel.Height = 11160
